$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb2"
$ws.Range("C2").Value = "Ephb2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 28.95628266666667
$ws.Range("H2").Value = 86.868848
$ws.Range("I2").Value = 0.5491054194301004
$ws.Range("J2").Value = 0.5491054194301005
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.882730333333333
$ws.Range("N2").Value = 26.648191
$ws.Range("O2").Value = 0.9229419957556332
$ws.Range("P2").Value = 0.9229419957556332
$ws.Range("Q2").Value = 257.2108503837742
$ws.Range("R2").Value = 2314.897653453968
$ws.Range("S2").Value = 0.5067924516890508
$ws.Range("T2").Value = 0.506792451689051

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb2"
$ws.Range("C3").Value = "Ephb2"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 28.95628266666667
$ws.Range("H3").Value = 86.868848
$ws.Range("I3").Value = 0.5491054194301004
$ws.Range("J3").Value = 0.5491054194301005
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.008116333333333333
$ws.Range("N3").Value = 0.024349
$ws.Range("O3").Value = 0.0008433110770878936
$ws.Range("P3").Value = 0.0008433110770878937
$ws.Range("Q3").Value = 0.2350188422168889
$ws.Range("R3").Value = 2.115169579952
$ws.Range("S3").Value = 0.0004630666826943975
$ws.Range("T3").Value = 0.0004630666826943977

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb2"
$ws.Range("C4").Value = "Ephb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 28.95628266666667
$ws.Range("H4").Value = 86.868848
$ws.Range("I4").Value = 0.5491054194301004
$ws.Range("J4").Value = 0.5491054194301005
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7335180000000001
$ws.Range("N4").Value = 2.200554
$ws.Range("O4").Value = 0.07621469316727886
$ws.Range("P4").Value = 0.07621469316727887
$ws.Range("Q4").Value = 21.239954549088
$ws.Range("R4").Value = 191.159590941792
$ws.Range("S4").Value = 0.04184990105835507
$ws.Range("T4").Value = 0.04184990105835508

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb2"
$ws.Range("C5").Value = "Ephb2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.691493
$ws.Range("H5").Value = 38.074479
$ws.Range("I5").Value = 0.2406720388519202
$ws.Range("J5").Value = 0.2406720388519202
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.882730333333333
$ws.Range("N5").Value = 26.648191
$ws.Range("O5").Value = 0.9229419957556332
$ws.Range("P5").Value = 0.9229419957556332
$ws.Range("Q5").Value = 112.7351098463877
$ws.Range("R5").Value = 1014.615988617489
$ws.Range("S5").Value = 0.2221263318605685
$ws.Range("T5").Value = 0.2221263318605685

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb2"
$ws.Range("C6").Value = "Ephb2"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.691493
$ws.Range("H6").Value = 38.074479
$ws.Range("I6").Value = 0.2406720388519202
$ws.Range("J6").Value = 0.2406720388519202
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.008116333333333333
$ws.Range("N6").Value = 0.024349
$ws.Range("O6").Value = 0.0008433110770878936
$ws.Range("P6").Value = 0.0008433110770878937
$ws.Range("Q6").Value = 0.1030083876856667
$ws.Range("R6").Value = 0.9270754891709999
$ws.Range("S6").Value = 0.0002029613963091522
$ws.Range("T6").Value = 0.0002029613963091522

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb2"
$ws.Range("C7").Value = "Ephb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.691493
$ws.Range("H7").Value = 38.074479
$ws.Range("I7").Value = 0.2406720388519202
$ws.Range("J7").Value = 0.2406720388519202
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7335180000000001
$ws.Range("N7").Value = 2.200554
$ws.Range("O7").Value = 0.07621469316727886
$ws.Range("P7").Value = 0.07621469316727887
$ws.Range("Q7").Value = 9.309438562374002
$ws.Range("R7").Value = 83.78494706136601
$ws.Range("S7").Value = 0.01834274559504251
$ws.Range("T7").Value = 0.01834274559504252

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Efnb2"
$ws.Range("C8").Value = "Ephb2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4888703333333334
$ws.Range("H8").Value = 1.466611
$ws.Range("I8").Value = 0.009270573592685367
$ws.Range("J8").Value = 0.009270573592685367
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.882730333333333
$ws.Range("N8").Value = 26.648191
$ws.Range("O8").Value = 0.9229419957556332
$ws.Range("P8").Value = 0.9229419957556332
$ws.Range("Q8").Value = 4.342503338966778
$ws.Range("R8").Value = 39.082530050701
$ws.Range("S8").Value = 0.008556201693432504
$ws.Range("T8").Value = 0.008556201693432504

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Efnb2"
$ws.Range("C9").Value = "Ephb2"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4888703333333334
$ws.Range("H9").Value = 1.466611
$ws.Range("I9").Value = 0.009270573592685367
$ws.Range("J9").Value = 0.009270573592685367
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.008116333333333333
$ws.Range("N9").Value = 0.024349
$ws.Range("O9").Value = 0.0008433110770878936
$ws.Range("P9").Value = 0.0008433110770878937
$ws.Range("Q9").Value = 0.003967834582111111
$ws.Range("R9").Value = 0.035710511239
$ws.Range("S9").Value = [double]"7.81797740167008E-06"
$ws.Range("T9").Value = [double]"7.817977401670082E-06"

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Efnb2"
$ws.Range("C10").Value = "Ephb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4888703333333334
$ws.Range("H10").Value = 1.466611
$ws.Range("I10").Value = 0.009270573592685367
$ws.Range("J10").Value = 0.009270573592685367
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7335180000000001
$ws.Range("N10").Value = 2.200554
$ws.Range("O10").Value = 0.07621469316727886
$ws.Range("P10").Value = 0.07621469316727887
$ws.Range("Q10").Value = 0.3585951891660001
$ws.Range("R10").Value = 3.227356702494001
$ws.Range("S10").Value = 0.0007065539218511932
$ws.Range("T10").Value = 0.0007065539218511933

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Efnb2"
$ws.Range("C11").Value = "Ephb2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.59691233333333
$ws.Range("H11").Value = 31.790737
$ws.Range("I11").Value = 0.200951968125294
$ws.Range("J11").Value = 0.200951968125294
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.882730333333333
$ws.Range("N11").Value = 26.648191
$ws.Range("O11").Value = 0.9229419957556332
$ws.Range("P11").Value = 0.9229419957556332
$ws.Range("Q11").Value = 94.12951462297411
$ws.Range("R11").Value = 847.1656316067669
$ws.Range("S11").Value = 0.1854670105125813
$ws.Range("T11").Value = 0.1854670105125813

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Efnb2"
$ws.Range("C12").Value = "Ephb2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 10.59691233333333
$ws.Range("H12").Value = 31.790737
$ws.Range("I12").Value = 0.200951968125294
$ws.Range("J12").Value = 0.200951968125294
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008116333333333333
$ws.Range("N12").Value = 0.024349
$ws.Range("O12").Value = 0.0008433110770878936
$ws.Range("P12").Value = 0.0008433110770878937
$ws.Range("Q12").Value = 0.08600807280144444
$ws.Range("R12").Value = 0.7740726552129999
$ws.Range("S12").Value = 0.0001694650206826738
$ws.Range("T12").Value = 0.0001694650206826738

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Efnb2"
$ws.Range("C13").Value = "Ephb2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 10.59691233333333
$ws.Range("H13").Value = 31.790737
$ws.Range("I13").Value = 0.200951968125294
$ws.Range("J13").Value = 0.200951968125294
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7335180000000001
$ws.Range("N13").Value = 2.200554
$ws.Range("O13").Value = 0.07621469316727886
$ws.Range("P13").Value = 0.07621469316727887
$ws.Range("Q13").Value = 7.773025940922001
$ws.Range("R13").Value = 69.95723346829801
$ws.Range("S13").Value = 0.01531549259203009
$ws.Range("T13").Value = 0.01531549259203009

